$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 860.25
$ws.Range("I80").Value = 293.45456
$ws.Range("K80").Value = 880.36368
$ws.Range("M80").Value = 117.63632
$ws.Range("H83").Value = 860.25
$ws.Range("I83").Value = 293.45456
$ws.Range("K83").Value = 2641.09104
$ws.Range("M83").Value = 2350.90896
$ws.Range("H113").Value = 3035.5293
$ws.Range("I113").Value = 2918
$ws.Range("J113").Value = 3251
$ws.Range("K113").Value = 2918
$ws.Range("L113").Value = 3251
$ws.Range("M113").Value = 336
$ws.Range("N113").Value = -9759
$ws.Range("H116").Value = 1885
$ws.Range("I116").Value = 1845
$ws.Range("J116").Value = 1925
$ws.Range("K116").Value = 1845
$ws.Range("L116").Value = 1925
$ws.Range("M116").Value = 1597
$ws.Range("N116").Value = -8809
$ws.Range("H121").Value = 1835
$ws.Range("J121").Value = 1835
$ws.Range("L121").Value = 5505
$ws.Range("N121").Value = -8999
$ws.Range("H131").Value = 1070.6666
$ws.Range("I131").Value = 764.8
$ws.Range("J131").Value = 2600
$ws.Range("K131").Value = 2294.4
$ws.Range("L131").Value = 7800
$ws.Range("M131").Value = 2745.6
$ws.Range("N131").Value = -17880
$ws.Range("H141").Value = 3157.8333
$ws.Range("I141").Value = 1761.875
$ws.Range("J141").Value = 5949.75
$ws.Range("K141").Value = 5285.625
$ws.Range("L141").Value = 17849.25
$ws.Range("M141").Value = -105.625
$ws.Range("N141").Value = -28209.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 186.66667
$ws.Range("I4").Value = 160
$ws.Range("J4").Value = 266.66666
$ws.Range("K4").Value = 160
$ws.Range("L4").Value = 266.66666
$ws.Range("M4").Value = -44
$ws.Range("N4").Value = -498.66666
$ws.Range("H5").Value = 55555800
$ws.Range("I5").Value = 55555800
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 55555800
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = -55555688
$ws.Range("N5").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("I76").Value = 29400
$ws.Range("J76").Value = 31833.334
$ws.Range("K76").Value = 29400
$ws.Range("L76").Value = 31833.334
$ws.Range("M76").Value = -29062
$ws.Range("N76").Value = -32509.334
$ws.Range("I79").Value = 29400
$ws.Range("J79").Value = 31833.334
$ws.Range("K79").Value = 29400
$ws.Range("L79").Value = 31833.334
$ws.Range("M79").Value = -28230
$ws.Range("N79").Value = -34173.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 55555800
$ws.Range("I4").Value = 55555800
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 55555800
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = -55555685
$ws.Range("N4").Value = 0
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 0
$ws.Range("L8").ClearContents()
$ws.Range("M8").Value = 1000
$ws.Range("N8").Value = -1280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 4558.4614
$ws.Range("J41").Value = 4558.4614
$ws.Range("L41").Value = 4558.4614
$ws.Range("N41").Value = -5414.4614
$ws.Range("H107").Value = 563.95654
$ws.Range("I107").Value = 513.4706
$ws.Range("K107").Value = 513.4706
$ws.Range("M107").Value = 1406.5294
$ws.Range("H132").Value = 2096.158
$ws.Range("I132").Value = 1427.6451
$ws.Range("J132").Value = 5056.7144
$ws.Range("K132").Value = 4282.9353
$ws.Range("L132").Value = 15170.1432
$ws.Range("M132").Value = -1752.9353
$ws.Range("N132").Value = -20230.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 372
$ws.Range("I2").Value = 441.89474
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 2651.36844
$ws.Range("L2").Value = 240
$ws.Range("M2").Value = -2538.36844
$ws.Range("N2").Value = -466

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4138340.2
$ws.Range("I11").Value = 5087712.5
$ws.Range("J11").Value = 2002253
$ws.Range("K11").Value = 5087712.5
$ws.Range("L11").Value = 2002253
$ws.Range("M11").Value = -5087573.5
$ws.Range("N11").Value = -2002531
$ws.Range("H12").Value = 5712500
$ws.Range("I12").Value = 6100000
$ws.Range("J12").Value = 3000000
$ws.Range("K12").Value = 6100000
$ws.Range("L12").Value = 3000000
$ws.Range("M12").Value = -6099860
$ws.Range("N12").Value = -3000280
$ws.Range("H14").Value = 39902.5
$ws.Range("I14").Value = 50000
$ws.Range("K14").Value = 50000
$ws.Range("M14").Value = -49832
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H80").Value = 6128.273
$ws.Range("I80").Value = 2381
$ws.Range("J80").Value = 9251
$ws.Range("K80").Value = 2381
$ws.Range("L80").Value = 9251
$ws.Range("M80").Value = -1383
$ws.Range("N80").Value = -11247
$ws.Range("H83").Value = 6128.273
$ws.Range("I83").Value = 2381
$ws.Range("J83").Value = 9251
$ws.Range("K83").Value = 11905
$ws.Range("L83").Value = 46255
$ws.Range("M83").Value = -6913
$ws.Range("N83").Value = -56239
$ws.Range("H122").Value = 43480428
$ws.Range("I122").Value = 90910650
$ws.Range("J122").Value = 2725.6667
$ws.Range("K122").Value = 272731950
$ws.Range("L122").Value = 8177.000100000001
$ws.Range("M122").Value = -272729500
$ws.Range("N122").Value = -13077.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 4009
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 5018
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 5018
$ws.Range("M30").Value = -2892
$ws.Range("N30").Value = -5234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 45428.75
$ws.Range("J133").Value = 45428.75
$ws.Range("L133").Value = 45428.75
$ws.Range("N133").Value = -55548.75
